$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.79589999999998
$ws.Range("D4").Value = -7.432099999999997
$ws.Range("A6").Value = -22.55000000000002
$ws.Range("A7").Value = -20.11429999999997
$ws.Range("D9").Value = -7.457699999999996
$ws.Range("D12").Value = -5.887999999999997
$ws.Range("A16").Value = -21.96570000000002
$ws.Range("D17").Value = -8.332799999999997
$ws.Range("D18").Value = -9.16279999999999
$ws.Range("D19").Value = -8.235299999999993
$ws.Range("A20").Value = -20.26329999999998
$ws.Range("D20").Value = -7.10769999999999
$ws.Range("D26").Value = -7.332000000000002
$ws.Range("A28").Value = -21.87769999999999
$ws.Range("A29").Value = -21.21099999999998
$ws.Range("D31").Value = -7.614499999999996
$ws.Range("A32").Value = -21.12949999999999
$ws.Range("D39").Value = -8.115599999999995
$ws.Range("A40").Value = -20.56670000000001
$ws.Range("D40").Value = -7.617099999999997
$ws.Range("D41").Value = -7.917399999999991
$ws.Range("D42").Value = -8.003899999999993
$ws.Range("D43").Value = -7.400400000000007
$ws.Range("A46").Value = -21.833
$ws.Range("D47").Value = -7.600900000000003
$ws.Range("D48").Value = -7.425899999999996
$ws.Range("A51").Value = -21.59489999999998
$ws.Range("A52").Value = -22.09599999999999
$ws.Range("A57").Value = -22.64470000000003
$ws.Range("A59").Value = -22.3171
$ws.Range("A62").Value = -22.1303
$ws.Range("D63").Value = -6.652999999999996
$ws.Range("D64").Value = -6.883699999999996
$ws.Range("A66").Value = -21.47989999999998
$ws.Range("A73").Value = -20.38279999999999
$ws.Range("A74").Value = -21.56499999999998
$ws.Range("D76").Value = -7.584399999999997
$ws.Range("D81").Value = -7.894900000000003
$ws.Range("D89").Value = -8.3446
$ws.Range("A92").Value = -21.60490000000001
$ws.Range("D94").Value = -5.9382
$ws.Range("A100").Value = -22.08780000000001

Write-Output "Applied 42 cell updates"
